$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Add the "Parameters" sheet first (so sheetId allocation matches the
#    original authoring order: Parameters gets sheetId 4, Sales data 2 gets
#    sheetId 5 below), even though Parameters ends up last in tab order.
#    Finish *everything* on this sheet right away: COM sheet handles/ranges
#    here resolve by tab position, so later inserts earlier in the tab order
#    (the "Sales data 2" copy) would otherwise make a stale reference point
#    at the wrong sheet.
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$parameters = $wb.Worksheets.Add($null, $lastSheet)
$parameters.Name = "Parameters"

$parameters.Range("B2").Value = "Param_Bloc1"

$parameters.Range("B3").Value = "Alpha"
$parameters.Range("C3").Value = "Beta"
$parameters.Range("D3").Value = "Gamma"

$parameters.Range("B4").Value = "aa"
$parameters.Range("C4").Value = "z"
$parameters.Range("D4").Value = 111

$parameters.Range("B5").Value = "bbb"
$parameters.Range("C5").Value = "y"
$parameters.Range("D5").Value = 222

$parameters.Range("B6").Value = "cccc"
$parameters.Range("C6").Value = "x"
$parameters.Range("D6").Value = 333

$parameters.Range("B7").Value = "ddddd"
$parameters.Range("C7").Value = "w"
$parameters.Range("D7").Value = 444

$parameters.Range("B11").Value = "Param_Block2"

$parameters.Range("B12").Value = "Alpha2"
$parameters.Range("C12").Value = "Beta2"
$parameters.Range("D12").Value = "Gamma2"

$parameters.Range("B13").Value = "aa"
$parameters.Range("C13").Value = "z"
$parameters.Range("D13").Value = 111

$parameters.Range("B14").Value = "bbb"
$parameters.Range("C14").Value = "y"
$parameters.Range("D14").Value = 222

$parameters.Range("B15").Value = "cccc"
$parameters.Range("C15").Value = "x"
$parameters.Range("D15").Value = 333

$parameters.Range("B16").Value = "ddddd"
$parameters.Range("C16").Value = "w"
$parameters.Range("D16").Value = 444

$parameters.Range("B12:D16").Select()

$wb.Names.Add('Block2', '=Parameters!$B$12:$D$16')
$wb.Names.Add('AlphaBlock', '=Parameters!$B$3:$D$7')

# ---------------------------------------------------------------------------
# 2) Duplicate "Sales data" into "Sales data 1" (renamed) + "Sales data 2"
#    (copy, keeps old VLOOKUP formulas), positioned right after the original.
# ---------------------------------------------------------------------------
$salesData1 = $wb.Worksheets.Item("Sales data")
$salesData1.Copy($null, $salesData1)

$salesData2 = $wb.Worksheets.Item(2)
$salesData2.Name = "Sales data 2"

$salesData1.Name = "Sales data 1"

# ---------------------------------------------------------------------------
# 3) On "Sales data 1", replace the VLOOKUP-driven G2:G8 with plain values
#    (the literal text "Country"), clearing the formula + the s="7" style.
# ---------------------------------------------------------------------------
$gRange = $salesData1.Range("G2:G8")
$gRange.ClearFormats()
$gRange.Value = "Country"

# ---------------------------------------------------------------------------
# 4) Selection bookkeeping: both "Sales data 1" and "Sales data 2" end up
#    with G1:G8 selected; "Sales data 1" stays the active tab.
# ---------------------------------------------------------------------------
$salesData2.Activate()
$salesData2.Range("G1:G8").Select()

$salesData1.Activate()
$salesData1.Range("G1:G8").Select()

# ---------------------------------------------------------------------------
# 5) Defined names local to each "Sales data N" sheet, pointing at its own
#    G1:G8.
# ---------------------------------------------------------------------------
$salesData1.Names.Add('SalesCountry', '=''Sales data 1''!$G$1:$G$8')
$salesData2.Names.Add('SalesCountry', '=''Sales data 2''!$G$1:$G$8')

# ---------------------------------------------------------------------------
# 6) Summary sheet: move selection from E39 to E23 (formulas auto-follow the
#    "Sales data" -> "Sales data 1" rename already performed above).
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Activate()
$summary.Range("E23").Select()

$salesData1.Activate()
